$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "564.92") need the
# NumberFormat forced to Text before assignment so Excel does not silently
# convert them to a numeric value (the source data keeps these as strings,
# e.g. "564.92", "7.07", "0.989"). Resetting the Style back to "Normal"
# afterwards avoids leaving a stray text-format style on the cell.

$ws.Range("D2").Value = '59.199.38'
$ws.Range("E2").Value = '  +0.86%  '

$ws.Range("D3").Value = '2.980.40'
$ws.Range("E3").Value = '  -0.60%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.41%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -0.24%  '

$ws.Range("D9").Value = '2.974.64'
$ws.Range("E9").Value = '  -0.45%  '

$ws.Range("E10").Value = '  +0.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.26'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.60%  '

$ws.Range("E12").Value = '  -1.17%  '

$ws.Range("E13").Value = '  -0.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.34%  '

$ws.Range("E15").Value = '  -0.45%  '

$ws.Range("D16").Value = '3.469.05'
$ws.Range("E16").Value = '  -0.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.16%  '

$ws.Range("D18").Value = '2.975.48'
$ws.Range("E18").Value = '  -0.53%  '

$ws.Range("D19").Value = '59.199.07'
$ws.Range("E19").Value = '  +1.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '434.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.66'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("E22").Value = '  +1.15%  '

$ws.Range("E23").Value = '  -1.39%  '

$ws.Range("E24").Value = '  -3.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.77%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.80%  '

$ws.Range("E28").Value = '  +0.11%  '

$ws.Range("E29").Value = '  +0.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.70'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.52%  '

$ws.Range("E32").Value = '  +0.60%  '

$ws.Range("E33").Value = '  +4.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.989'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.88'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.65%  '

$ws.Range("D36").Value = '0.0₃0760'
$ws.Range("E36").Value = '  +2.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.04'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.79'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '393.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0351'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.42%  '

$ws.Range("D43").Value = '2.701.41'
$ws.Range("E43").Value = '  -1.18%  '

$ws.Range("E44").Value = '  -3.64%  '

$ws.Range("E45").Value = '  +0.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.79'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.21%  '

$ws.Range("E49").Value = '  -0.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.40%  '
